{"js": "// Remove the lone \"\\n\" paragraph and the \"Please read Using Pressure\n// Canners...\" paragraph that used to sit between the \"Quality:\" paragraph\n// and the \"\\n\" paragraph preceding \"Procedure:\".\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nconst targets = [];\nfor (let i = 0; i < items.length; i++) {\n  const text = items[i].text;\n  if (text.indexOf(\"Please read Using Pressure Canners before beginning\") !== -1) {\n    targets.push(items[i]);\n    // The paragraph immediately before it is the standalone \"\\n\" paragraph\n    // that is also being removed.\n    if (i - 1 >= 0 && items[i - 1].text === \"\\\\n\") {\n      targets.push(items[i - 1]);\n    }\n    break;\n  }\n}\n\nfor (const paragraph of targets) {\n  paragraph.delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the lone \"\\n\" paragraph and the \"Please read Using Pressure\n# Canners...\" paragraph that used to sit between the \"Quality:\" paragraph\n# and the \"\\n\" paragraph preceding \"Procedure:\".\n$d = $word.ActiveDocument\n\n$wdParagraph = 4\n\n$rng = $d.Content\n$found = $rng.Find.Execute(\"Please read Using Pressure Canners before beginning\")\nif ($found) {\n    # Grow the found range to cover its whole paragraph (text + mark).\n    $rng.Expand($wdParagraph) | Out-Null\n\n    # The paragraph right before it is the standalone \"\\n\" paragraph that\n    # also needs to go. Pull its start in so one Delete removes both.\n    # (Range.Text carries a trailing paragraph-mark control character, so\n    # trim it off before comparing against the literal \"\\n\" run text.)\n    $prev = $rng.Previous($wdParagraph, 1)\n    $prevText = $prev.Text.TrimEnd([char]13, [char]7)\n    if ($prev -ne $null -and $prevText -eq \"\\n\") {\n        $rng.Start = $prev.Start\n    }\n\n    $rng.Delete()\n}\n"}
